$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.162.58"
$ws.Range("E2").Value = "  +1.00%  "

# Row 3
$ws.Range("D3").Value = "2.662.08"
$ws.Range("E3").Value = "  +1.97%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'532.38"
$ws.Range("E5").Value = "  +4.07%  "

# Row 6
$ws.Range("D6").Value = "'156.84"
$ws.Range("E6").Value = "  +1.69%  "

# Row 7
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  +0.63%  "

# Row 9
$ws.Range("E9").Value = "  -2.34%  "

# Row 10
$ws.Range("E10").Value = "  +5.06%  "

# Row 11
$ws.Range("E11").Value = "  +2.27%  "

# Row 13
$ws.Range("D13").Value = "3.128.27"
$ws.Range("E13").Value = "  +1.96%  "

# Row 14
$ws.Range("D14").Value = "61.144.93"
$ws.Range("E14").Value = "  +1.06%  "

# Row 15
$ws.Range("D15").Value = "'22.10"
$ws.Range("E15").Value = "  +2.42%  "

# Row 16
$ws.Range("E16").Value = "  +2.28%  "

# Row 17
$ws.Range("D17").Value = "2.673.59"
$ws.Range("E17").Value = "  +2.00%  "

# Row 18
$ws.Range("D18").Value = "'4.79"
$ws.Range("E18").Value = "  +0.78%  "

# Row 19
$ws.Range("D19").Value = "'356.20"
$ws.Range("E19").Value = "  +0.85%  "

# Row 20
$ws.Range("D20").Value = "'10.73"
$ws.Range("E20").Value = "  +1.49%  "

# Row 21
$ws.Range("D21").Value = "'6.31"
$ws.Range("E21").Value = "  +2.18%  "

# Row 22
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
$ws.Range("D23").Value = "'61.61"
$ws.Range("E23").Value = "  +1.57%  "

# Row 24
$ws.Range("E24").Value = "  +2.54%  "

# Row 25
$ws.Range("E25").Value = "  +1.91%  "

# Row 27
$ws.Range("D27").Value = "0.0₃0865"
$ws.Range("E27").Value = "  +2.94%  "

# Row 28
$ws.Range("D28").Value = "'7.44"
$ws.Range("E28").Value = "  +1.44%  "

# Row 30
$ws.Range("D30").Value = "'6.19"
$ws.Range("E30").Value = "  +6.72%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.64"
$ws.Range("E31").Value = "  +4.22%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'19.60"
$ws.Range("E32").Value = "  +1.06%  "

# Row 33
$ws.Range("D33").Value = "'149.99"
$ws.Range("E33").Value = "  -0.89%  "

# Row 34
$ws.Range("E34").Value = "  +4.43%  "

# Row 35
$ws.Range("E35").Value = "  +1.38%  "

# Row 36
$ws.Range("D36").Value = "'0.916"
$ws.Range("E36").Value = "  +8.82%  "

# Row 37
$ws.Range("D37").Value = "'0.886"
$ws.Range("E37").Value = "  -0.29%  "

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'309.62"
$ws.Range("E38").Value = "  +6.47%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.51"
$ws.Range("E39").Value = "  +1.31%  "

# Row 40
$ws.Range("E40").Value = "  +1.79%  "

# Row 41
$ws.Range("D41").Value = "'0.651"
$ws.Range("E41").Value = "  +4.07%  "

# Row 42
$ws.Range("E42").Value = "  +1.45%  "

# Row 43
$ws.Range("D43").Value = "'20.51"
$ws.Range("E43").Value = "  +3.83%  "

# Row 44
$ws.Range("E44").Value = "  +2.38%  "

# Row 45
$ws.Range("E45").Value = "  +0.07%  "

# Row 46
$ws.Range("D46").Value = "'5.03"
$ws.Range("E46").Value = "  +2.61%  "

# Row 47
$ws.Range("D47").Value = "'0.0241"
$ws.Range("E47").Value = "  +2.83%  "

# Row 48
$ws.Range("D48").Value = "'10.35"
$ws.Range("E48").Value = "  +0.39%  "

# Row 49
$ws.Range("D49").Value = "'19.14"
$ws.Range("E49").Value = "  +8.88%  "

# Row 50
$ws.Range("D50").Value = "2.001.58"
$ws.Range("E50").Value = "  +0.19%  "

# Row 51
$ws.Range("E51").Value = "  +3.22%  "
